$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows for row 2 and row 3 got swapped for columns
# A (Id), Q (Ost), R (Nord), Z (Starttid) and AB (Sluttid).
# Capture the "before" values first, then write them back swapped.

$cols = @("A", "Q", "R", "Z", "AB")

$row2vals = @{}
$row3vals = @{}

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $row2vals[$col] = $ws.Range($addr2).Value()
    $row3vals[$col] = $ws.Range($addr3).Value()
}

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $ws.Range($addr2).Value = $row3vals[$col]
    $ws.Range($addr3).Value = $row2vals[$col]
}
